$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-33 for columns D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo), M (Precio promedio ponderado), P (Precio $/Kg)
$data = @{
    2 = @{ D = 44509; J = 100; K = 15000; L = 16000; M = 15500; P = 1192 }
    3 = @{ D = 44610; J = 50; K = 17000; L = 18000; M = 17400; P = 1338 }
    4 = @{ D = 44159; J = 60; K = 30000; L = 32000; M = 31000; P = 2385 }
    5 = @{ D = 44433; J = 100; K = 13000; L = 14000; M = 13500; P = 1038 }
    6 = @{ D = 44316; J = 50; K = 27000; L = 28000; M = 27400; P = 2108 }
    7 = @{ D = 44350; J = 40; K = 23000; L = 25000; M = 24000; P = 1846 }
    8 = @{ D = 44503; J = 35; K = 15000; L = 16000; M = 15429; P = 1187 }
    9 = @{ D = 44708; J = 50; K = 13000; L = 14000; M = 13600; P = 1046 }
    10 = @{ D = 44383; J = 50; K = 15000; L = 16000; M = 15400; P = 1185 }
    11 = @{ D = 44474; J = 40; K = 13000; L = 14000; M = 13500; P = 1038 }
    12 = @{ D = 44308; J = 50; K = 26000; L = 27000; M = 26400; P = 2031 }
    13 = @{ D = 44488; J = 40; K = 16000; L = 17000; M = 16500; P = 1269 }
    14 = @{ D = 44313; J = 50; K = 25000; L = 26000; M = 25600; P = 1969 }
    15 = @{ D = 44377; J = 40; K = 14000; L = 15000; M = 14500; P = 1115 }
    16 = @{ D = 44467; J = 100; K = 13000; L = 14000; M = 13500; P = 1038 }
    17 = @{ D = 44741; J = 50; K = 14000; L = 15000; M = 14400; P = 1108 }
    18 = @{ D = 44523; J = 40; K = 15000; L = 16000; M = 15500; P = 1192 }
    19 = @{ D = 44691; J = 100; K = 12000; L = 13000; M = 12500; P = 962 }
    20 = @{ D = 44510; J = 40; K = 15000; L = 16000; M = 15500; P = 1192 }
    21 = @{ D = 44435; J = 100; K = 13000; L = 14000; M = 13500; P = 1038 }
    22 = @{ D = 44327; J = 50; K = 24000; L = 25000; M = 24400; P = 1877 }
    23 = @{ D = 44462; J = 60; K = 14000; L = 15000; M = 14500; P = 1115 }
    24 = @{ D = 44453; J = 50; K = 14000; L = 15000; M = 14600; P = 1123 }
    25 = @{ D = 44334; J = 50; K = 26000; L = 28000; M = 27200; P = 2092 }
    26 = @{ D = 44425; J = 60; K = 14000; L = 15000; M = 14500; P = 1115 }
    27 = @{ D = 44719; J = 50; K = 13000; L = 14000; M = 13400; P = 1031 }
    28 = @{ D = 44705; J = 50; K = 10000; L = 11000; M = 10400; P = 800 }
    29 = @{ D = 44664; J = 50; K = 11000; L = 12000; M = 11600; P = 892 }
    30 = @{ D = 44355; J = 60; K = 18000; L = 20000; M = 19000; P = 1462 }
    31 = @{ D = 44362; J = 40; K = 15000; L = 16000; M = 15500; P = 1192 }
    32 = @{ D = 44320; J = 50; K = 26000; L = 28000; M = 26800; P = 2062 }
    33 = @{ D = 44264; J = 40; K = 30000; L = 32000; M = 31000; P = 2385 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item([int]$row, 4).Value = $vals.D
    $ws.Cells.Item([int]$row, 10).Value = $vals.J
    $ws.Cells.Item([int]$row, 11).Value = $vals.K
    $ws.Cells.Item([int]$row, 12).Value = $vals.L
    $ws.Cells.Item([int]$row, 13).Value = $vals.M
    $ws.Cells.Item([int]$row, 16).Value = $vals.P
}

Write-Output "Done updating rows"